# Applies the cryptos-list price/volume refresh described in the commit
# message "Updated cryptos list on Sun Jun 16 23:35:36 UTC 2024 with GitHub
# Actions" -- 76 cell updates across rows 2-51 of the active sheet, including
# a reordering of the Kaspa / RenzoRestakedETH rows (31 <-> 32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.649.36"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "3.618.35"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.49"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.35"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("D7").Value = "3.616.90"
$ws.Range("E7").Value = "  +1.18%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.416"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "4.232.16"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "3.578.46"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "66.741.75"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.66"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.36"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +4.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.44"
$ws.Range("E27").Value = "  +6.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.52"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.160"
$ws.Range("E31").Value = "  +4.25%  "
$ws.Range("B32").Value = "RenzoRestakedETH"
$ws.Range("C32").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D32").Value = "3.615.71"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.48"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.67"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.29"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.902"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "46.26"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.61"
$ws.Range("E45").Value = "  +8.43%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.10"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.97"
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("E51").Value = "  +2.21%  "
